# Adjust ML Labels: flip the "Label" (column C) value from 1 to 0 for the
# rows where the classifier's prediction changed, then recalc the dependent
# formulas (D, E, H, I, J columns) automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$rowsToZero = @(3, 8, 10, 11, 12, 13, 15, 19, 20, 21, 23, 25, 28, 32, 38, 45, 47)

foreach ($r in $rowsToZero) {
    $ws.Cells.Item($r, 3).Value = 0
}

$excel.Calculate()

# Leave the cursor where the author left it when done reviewing the results.
$ws.Range("J3").Select()

$wb.Save()
